$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.547.11"
$ws.Range("E2").Value = "  -2.63%  "
$ws.Range("D3").Value = "3.367.54"
$ws.Range("E3").Value = "  -4.45%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'555.94"
$ws.Range("E5").Value = "  -4.72%  "
$ws.Range("D6").Value = "'176.33"
$ws.Range("E6").Value = "  -1.77%  "
$ws.Range("D7").Value = "'0.617"
$ws.Range("E7").Value = "  -2.38%  "
$ws.Range("D8").Value = "3.362.38"
$ws.Range("E8").Value = "  -4.35%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  -2.00%  "
$ws.Range("D11").Value = "'0.161"
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("D12").Value = "'54.91"
$ws.Range("E12").Value = "  -2.07%  "
$ws.Range("E13").Value = "  -3.04%  "
$ws.Range("E14").Value = "  -2.65%  "
$ws.Range("D15").Value = "3.907.62"
$ws.Range("E15").Value = "  -4.26%  "
$ws.Range("D16").Value = "'18.41"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "3.374.86"
$ws.Range("E17").Value = "  -4.27%  "
$ws.Range("D19").Value = "'11.84"
$ws.Range("E19").Value = "  -2.36%  "
$ws.Range("D20").Value = "64.476.14"
$ws.Range("E20").Value = "  -2.60%  "
$ws.Range("E21").Value = "  -2.99%  "
$ws.Range("D22").Value = "'431.75"
$ws.Range("E22").Value = "  +3.52%  "
$ws.Range("D24").Value = "'4.11"
$ws.Range("E24").Value = "  -4.92%  "
$ws.Range("D25").Value = "'84.24"
$ws.Range("E26").Value = "  -1.39%  "
$ws.Range("D27").Value = "'10.81"
$ws.Range("E27").Value = "  -2.80%  "
$ws.Range("E28").Value = "  -1.11%  "
$ws.Range("D29").Value = "'8.78"
$ws.Range("E29").Value = "  -4.06%  "
$ws.Range("D30").Value = "'29.75"
$ws.Range("E30").Value = "  -2.39%  "
$ws.Range("D31").Value = "'6.67"
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("D32").Value = "'11.50"
$ws.Range("E32").Value = "  -2.43%  "
$ws.Range("D33").Value = "'573.46"
$ws.Range("E33").Value = "  -5.70%  "
$ws.Range("E34").Value = "  -3.08%  "
$ws.Range("D35").Value = "'58.54"
$ws.Range("E35").Value = "  -2.36%  "
$ws.Range("D37").Value = "'0.142"
$ws.Range("E37").Value = "  -7.73%  "
$ws.Range("D38").Value = "'3.52"
$ws.Range("E38").Value = "  -3.68%  "
$ws.Range("D39").Value = "'35.77"
$ws.Range("E39").Value = "  -4.20%  "
$ws.Range("D40").Value = "0.0₃0758"
$ws.Range("E40").Value = "  -5.74%  "
$ws.Range("E41").Value = "  -4.07%  "
$ws.Range("D42").Value = "3.115.21"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("E44").Value = "  -6.08%  "
$ws.Range("D45").Value = "'3.27"
$ws.Range("E45").Value = "  -3.40%  "
$ws.Range("E46").Value = "  -2.91%  "
$ws.Range("E47").Value = "  -3.67%  "
$ws.Range("E48").Value = "  -2.16%  "
$ws.Range("E49").Value = "  -3.69%  "
$ws.Range("D50").Value = "'134.88"
$ws.Range("E50").Value = "  -3.58%  "
$ws.Range("D51").Value = "'8.26"
$ws.Range("E51").Value = "  -5.14%  "
